# "new change 16 june"
# - Sheet1!B8 password value changes from "Fosroc@1" to "Fosroc@3"
# - SEBS_Devloper!A3 contact name changes from "Raj Kumar" to "ravi varma"

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B8").Value = "Fosroc@3"
[void]$ws1.Range("B8").Select()

$ws2 = $wb.Worksheets.Item("SEBS_Devloper")
$ws2.Range("A3").Value = "ravi varma"
